$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "TC_PM_COA_SEC_ProgramCodeProfile_AddNewProfile"
$ws.Range("B10").Value = "TC_PM_COA_SEC_ProgramCodeProfile_EditProfile"
$ws.Range("B13").Value = "TC_PM_COA_SEC_ProgramCodeProfile_AddNewProfile"

$ws.Columns.Item(2).ColumnWidth = 46.5

$ws.Range("B7").Select()
